# Applies the "cryptos list" data refresh described by the commit
# "Updated cryptos list on Fri Sep 15 07:54:16 UTC 2023 with GitHub Actions".
#
# For cells whose new text happens to look like a plain number (e.g. "212.60",
# "0.497") a leading apostrophe is used so Excel stores the value as text
# (matching the workbook's original inlineStr/text representation) instead of
# silently parsing it into a numeric cell and mangling the formatting
# (trailing zeros, float rounding, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "26.641.02"
$ws.Range("E2").Value = "  +1.01%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.631.95"
$ws.Range("E3").Value = "  +0.56%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.21%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "'212.60"
$ws.Range("E5").Value = "  +0.09%  "

# --- Row 6 (XRP) ---
$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "  +2.72%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.20%  "

# --- Row 8 (Cardano) ---
$ws.Range("E8").Value = "  +1.59%  "

# --- Row 9 (Dogecoin) ---
$ws.Range("E9").Value = "  +0.63%  "

# --- Row 10 (Solana) ---
$ws.Range("E10").Value = "  +1.27%  "

# --- Row 11 (TRON) ---
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +3.55%  "

# --- Row 12 (WrappedliquidstakedEther2.0) ---
$ws.Range("D12").Value = "1.859.82"
$ws.Range("E12").Value = "  +0.62%  "

# --- Row 13: was Polkadot, now WrappedEther ---
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.661.36"
$ws.Range("E13").Value = "  +2.42%  "

# --- Row 14: was WrappedEther, now Polkadot ---
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.08"
$ws.Range("E14").Value = "  +1.38%  "

# --- Row 15 (Polygon) ---
$ws.Range("E15").Value = "  +0.65%  "

# --- Row 16 (WrappedBTC) ---
$ws.Range("D16").Value = "26.648.65"
$ws.Range("E16").Value = "  +1.08%  "

# --- Row 17 (Litecoin) ---
$ws.Range("D17").Value = "'63.18"
$ws.Range("E17").Value = "  +0.89%  "

# --- Row 18 (ShibaInu) ---
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.61%  "

# --- Row 19 (BitcoinCash) ---
$ws.Range("D19").Value = "'217.90"
$ws.Range("E19").Value = "  +7.66%  "

# --- Row 20 (Dai) ---
$ws.Range("E20").Value = "  +0.11%  "

# --- Row 21 (Uniswap) ---
$ws.Range("D21").Value = "'4.27"
$ws.Range("E21").Value = "  -0.21%  "

# --- Row 22 (Chainlink) ---
$ws.Range("E22").Value = "  +2.01%  "

# --- Row 23 (Avalanche) ---
$ws.Range("D23").Value = "'9.36"
$ws.Range("E23").Value = "  -0.01%  "

# --- Row 24 (Toncoin) ---
$ws.Range("E24").Value = "  +3.33%  "

# --- Row 25 (Monero) ---
$ws.Range("D25").Value = "'148.05"
$ws.Range("E25").Value = "  +2.51%  "

# --- Row 26 (BinanceUSD) ---
$ws.Range("E26").Value = "  +0.18%  "

# --- Row 27 (Stellar) ---
$ws.Range("E27").Value = "  +0.73%  "

# --- Row 28 (Cosmos) ---
$ws.Range("D28").Value = "'6.84"
$ws.Range("E28").Value = "  +4.03%  "

# --- Row 29 (EthereumClassic) ---
$ws.Range("D29").Value = "'15.48"
$ws.Range("E29").Value = "  +1.88%  "

# --- Row 30 (Hedera) ---
$ws.Range("E30").Value = "  -2.27%  "

# --- Row 31 (PancakeSwap) ---
$ws.Range("E31").Value = "  -0.05%  "

# --- Row 32 (Filecoin) ---
$ws.Range("E32").Value = "  +3.25%  "

# --- Row 33 (InternetComputer(DFINITY)) ---
$ws.Range("E33").Value = "  +1.40%  "

# --- Row 34 (LidoDAOToken) ---
$ws.Range("E34").Value = "  +0.14%  "

# --- Row 35 (HuobiToken) ---
$ws.Range("E35").Value = "  +0.70%  "

# --- Row 36 (Maker) ---
$ws.Range("D36").Value = "1.209.01"
$ws.Range("E36").Value = "  +2.72%  "

# --- Row 37 (VeChain) ---
$ws.Range("E37").Value = "  +4.99%  "

# --- Row 38 (ARBITRUM) ---
$ws.Range("E38").Value = "  -0.44%  "

# --- Row 39 (PaxDollar) ---
$ws.Range("E39").Value = "  +0.28%  "

# --- Row 40 (ImmutableX) ---
$ws.Range("D40").Value = "'0.499"
$ws.Range("E40").Value = "  +0.32%  "

# --- Row 41 (MXToken) ---
$ws.Range("E41").Value = "  -1.96%  "

# --- Row 42 (FraxShare) ---
$ws.Range("D42").Value = "'5.40"
$ws.Range("E42").Value = "  +0.41%  "

# --- Row 43 (TrustWalletToken) ---
$ws.Range("D43").Value = "'0.791"
$ws.Range("E43").Value = "  +0.44%  "

# --- Row 44 (RocketPoolETH) ---
$ws.Range("D44").Value = "1.775.33"
$ws.Range("E44").Value = "  +0.83%  "

# --- Row 45 (Quant) ---
$ws.Range("D45").Value = "'92.78"
$ws.Range("E45").Value = "  +0.09%  "

# --- Row 46 (RenderToken) ---
$ws.Range("E46").Value = "  +0.63%  "

# --- Row 47 (Aave) ---
$ws.Range("D47").Value = "'54.64"
$ws.Range("E47").Value = "  +1.41%  "

# --- Row 48 (Cronos) ---
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  +0.91%  "

# --- Row 49 (EnergySwap) ---
$ws.Range("D49").Value = "'7.60"
$ws.Range("E49").Value = "  +3.86%  "

# --- Row 50 (Mantle) ---
$ws.Range("E50").Value = "  -0.04%  "

# --- Row 51 (USDD) ---
$ws.Range("E51").Value = "  +0.30%  "
